$d = $word.ActiveDocument

# The document contains 7 paragraphs shaped like:
#   <w:r><id></w:r><w:r>p063v_N</w:r><w:r></id></w:r>
# (three separate runs with differing rPr), each forming the whole text of
# its own paragraph ("<id>p063v_N</id>"). The edit merges those three runs
# into a single run (keeping the formatting/rPr of the first "<id>" run)
# whose text is the concatenation "<id>p063v_N</id>".
#
# We collect the matching paragraph ranges first (by their exact text),
# then process them in reverse document order so that earlier offsets are
# not invalidated by edits made later in the document.

$targets = New-Object System.Collections.ArrayList

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^<id>p063v_\d+</id>\r?$") {
        $r = $p.Range
        [void]$targets.Add(@($r.Start, $r.End))
    }
}

for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $bounds = $targets[$i]
    $start = $bounds[0]
    $end = $bounds[1] - 1   # exclude the paragraph mark

    $full = $d.Range($start, $end)
    $fullText = $full.Text   # e.g. "<id>p063v_3</id>"

    $openLen = 4             # length of "<id>"
    $closeLen = 5            # length of "</id>"
    $middle = $fullText.Substring($openLen, $fullText.Length - $openLen - $closeLen)

    $openEnd = $start + $openLen

    # Delete the third run ("</id>") and the second run (the id value),
    # from the back of the range forward, leaving only the first run.
    $closeRange = $d.Range($start + $fullText.Length - $closeLen, $end)
    $closeRange.Delete()

    $middleRange = $d.Range($openEnd, $start + $fullText.Length - $closeLen)
    $middleRange.Delete()

    # Re-append the removed text onto the end of the first run so it is
    # absorbed into that run instead of creating a brand-new run.
    $firstRange = $d.Range($start, $openEnd)
    $firstRange.InsertAfter($middle + "</id>")
}
